# Rename the inline picture placeholders in the document's headers/footers.
#
# The BTec logo (stored in the "first page" header) goes from image1.jpg to
# image2.jpg, and the two Pearson logos (stored in the default + first-page
# footers) both go from image2.png to image1.png. Only the `name=` shown by
# InlineShapes (the OOXML docPr/name) changes - the alt text (descr) stays
# untouched.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Header (first page) : BTec_Logo-Orange : image1.jpg -> image2.jpg ---
$hdr = $sec.Headers.Item(2)
if ($hdr.Exists -and $hdr.Range.InlineShapes.Count -ge 1) {
    $shape = $hdr.Range.InlineShapes.Item(1)
    # Re-fetch the shape through its own Range so the rename reliably
    # writes back to the underlying header/footer part.
    $shape.Range.InlineShapes.Item(1).Name = "image2.jpg"
}

# --- Footer (default) : PearsonLogo : image2.png -> image1.png ---
$ftrDefault = $sec.Footers.Item(1)
if ($ftrDefault.Exists -and $ftrDefault.Range.InlineShapes.Count -ge 1) {
    $shape = $ftrDefault.Range.InlineShapes.Item(1)
    $shape.Range.InlineShapes.Item(1).Name = "image1.png"
}

# --- Footer (first page) : PearsonLogo : image2.png -> image1.png ---
$ftrFirst = $sec.Footers.Item(2)
if ($ftrFirst.Exists -and $ftrFirst.Range.InlineShapes.Count -ge 1) {
    $shape = $ftrFirst.Range.InlineShapes.Item(1)
    $shape.Range.InlineShapes.Item(1).Name = "image1.png"
}

Write-Output "Renamed header/footer logo inline shapes."
